$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "throttle fully open" value (C4): 26000 -> 24130
$ws.Range("C4").Value = 24130

# Row 8: B8 value changes from 5 -> 4 (formula in C8 stays the same, recalculates)
$ws.Range("B8").Value = 4

# Row 9: B9 value changes from 1.8 -> 5 (formula in C9 stays the same, recalculates)
$ws.Range("B9").Value = 5

# New row 10: B10 = 1.8, C10 = formula copied down from C9
$ws.Range("B10").Value = 1.8
$ws.Range("C10").Formula = "=B10/`$B`$4*`$C`$5+`$C`$3"
$ws.Range("C10").NumberFormat = $ws.Range("C9").NumberFormat

# New row 11: B11 = 100.8, C11 = formula copied down from C9
$ws.Range("B11").Value = 100.8
$ws.Range("C11").Formula = "=B11/`$B`$4*`$C`$5+`$C`$3"
$ws.Range("C11").NumberFormat = $ws.Range("C9").NumberFormat

# Update the active selection to match the target (E9)
$ws.Range("E9").Select()

$excel.Calculate()
